$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we touch to remain Text, matching the
# workbook author convention of storing prices as literal strings (not
# numbers) even when they look numeric (e.g. "215.31").
$priceCells = @(
    "D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12",
    "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22",
    "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32",
    "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42",
    "D43", "D44", "D45", "D47", "D48", "D49", "D50"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "25.901.57"
$ws.Range("E2").Value = "  -0.91%  "

# Row 3
$ws.Range("D3").Value = "1.638.16"
$ws.Range("E3").Value = "  -0.73%  "

# Row 4
$ws.Range("E4").Value = "  -0.28%  "

# Row 5
$ws.Range("D5").Value = "215.31"
$ws.Range("E5").Value = "  -0.02%  "

# Row 6
$ws.Range("D6").Value = "0.5031"
$ws.Range("E6").Value = "  -1.65%  "

# Row 7
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.30%  "

# Row 8
$ws.Range("D8").Value = "0.2564"
$ws.Range("E8").Value = "  -1.02%  "

# Row 9
$ws.Range("D9").Value = "0.06389"
$ws.Range("E9").Value = "  -0.59%  "

# Row 10
$ws.Range("D10").Value = "19.67"
$ws.Range("E10").Value = "  -1.14%  "

# Row 11
$ws.Range("D11").Value = "0.07727"
$ws.Range("E11").Value = "  -0.71%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.659.43"
$ws.Range("E12").Value = "  +0.40%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.261"
$ws.Range("E13").Value = "  -0.57%  "

# Row 14
$ws.Range("D14").Value = "1.864.70"
$ws.Range("E14").Value = "  -0.69%  "

# Row 15
$ws.Range("D15").Value = "0.5446"
$ws.Range("E15").Value = "  -1.10%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7908"
$ws.Range("E16").Value = "  -1.14%  "

# Row 17
$ws.Range("D17").Value = "64.32"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18
$ws.Range("D18").Value = "25.913.50"
$ws.Range("E18").Value = "  -0.92%  "

# Row 19
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.23%  "

# Row 20
$ws.Range("D20").Value = "203.12"
$ws.Range("E20").Value = "  -3.28%  "

# Row 21
$ws.Range("D21").Value = "4.378"
$ws.Range("E21").Value = "  -0.30%  "

# Row 22
$ws.Range("D22").Value = "9.905"
$ws.Range("E22").Value = "  -1.40%  "

# Row 23
$ws.Range("D23").Value = "5.986"
$ws.Range("E23").Value = "  -0.82%  "

# Row 24
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  -0.17%  "

# Row 25
$ws.Range("D25").Value = "1.930"
$ws.Range("E25").Value = "  +10.61%  "

# Row 26
$ws.Range("D26").Value = "141.23"
$ws.Range("E26").Value = "  -1.86%  "

# Row 27
$ws.Range("D27").Value = "0.1136"
$ws.Range("E27").Value = "  -3.27%  "

# Row 28
$ws.Range("D28").Value = "15.67"
$ws.Range("E28").Value = "  -0.89%  "

# Row 29
$ws.Range("D29").Value = "6.728"
$ws.Range("E29").Value = "  -3.33%  "

# Row 30
$ws.Range("D30").Value = "1.243"
$ws.Range("E30").Value = "  +0.20%  "

# Row 31
$ws.Range("D31").Value = "0.04937"
$ws.Range("E31").Value = "  -3.22%  "

# Row 32
$ws.Range("D32").Value = "3.277"
$ws.Range("E32").Value = "  -2.09%  "

# Row 33
$ws.Range("D33").Value = "3.183"
$ws.Range("E33").Value = "  -0.87%  "

# Row 34
$ws.Range("D34").Value = "1.544"
$ws.Range("E34").Value = "  -0.80%  "

# Row 35
$ws.Range("D35").Value = "2.376"
$ws.Range("E35").Value = "  +1.07%  "

# Row 36
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "2.630"
$ws.Range("E36").Value = "  -3.94%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "0.8942"
$ws.Range("E37").Value = "  -3.16%  "

# Row 38
$ws.Range("D38").Value = "1.161.05"
$ws.Range("E38").Value = "  +0.02%  "

# Row 39
$ws.Range("D39").Value = "0.5602"
$ws.Range("E39").Value = "  -1.86%  "

# Row 40
$ws.Range("D40").Value = "0.01564"
$ws.Range("E40").Value = "  -1.28%  "

# Row 41
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  -0.22%  "

# Row 42
$ws.Range("D42").Value = "5.718"
$ws.Range("E42").Value = "  +1.25%  "

# Row 43
$ws.Range("D43").Value = "0.8078"
$ws.Range("E43").Value = "  -2.12%  "

# Row 44
$ws.Range("D44").Value = "99.78"
$ws.Range("E44").Value = "  -0.30%  "

# Row 45
$ws.Range("D45").Value = "1.776.31"
$ws.Range("E45").Value = "  -0.69%  "

# Row 46
$ws.Range("E46").Value = "  -0.05%  "

# Row 47
$ws.Range("D47").Value = "0.4524"
$ws.Range("E47").Value = "  -0.59%  "

# Row 48
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.10%  "

# Row 49
$ws.Range("D49").Value = "55.00"
$ws.Range("E49").Value = "  -0.74%  "

# Row 50
$ws.Range("D50").Value = "0.05055"
$ws.Range("E50").Value = "  -0.46%  "

# Row 51
$ws.Range("E51").Value = "  -0.40%  "
